$wb = $excel.ActiveWorkbook

# Mapping of worksheet name -> cell:value updates (column F "想去人数")
$updates = @{
    "展览" = @{
        "F3"  = 1353
        "F5"  = 105
        "F7"  = 11664
        "F8"  = 4383
        "F9"  = 30
        "F14" = 1094
        "F16" = 39
        "F17" = 5082
        "F21" = 11331
    }
    "全部类型" = @{
        "F3"  = 1353
        "F5"  = 105
        "F7"  = 11664
        "F8"  = 4383
        "F9"  = 30
        "F15" = 1094
        "F17" = 39
        "F18" = 5082
        "F22" = 11331
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
